# Applies cryptos list update (Sun Mar 10 18:31:11 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.375.32"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "3.911.43"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Formula = "'530.12"
$ws.Range("E5").Value = "  +9.50%  "
$ws.Range("D6").Formula = "'144.09"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Formula = "'0.717"
$ws.Range("E9").Value = "  -3.17%  "
$ws.Range("E10").Value = "  -4.36%  "
$ws.Range("E11").Value = "  -5.12%  "
$ws.Range("D12").Formula = "'42.06"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").Value = "4.538.25"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Formula = "'10.27"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").Value = "3.914.78"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("E16").Value = "  +8.87%  "
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Formula = "'19.72"
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("D20").Value = "69.319.94"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").Formula = "'427.80"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("E22").Value = "  -5.43%  "
$ws.Range("D23").Formula = "'88.63"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("E24").Value = "  -4.71%  "
$ws.Range("D25").Formula = "'4.07"
$ws.Range("E25").Value = "  +10.15%  "
$ws.Range("D26").Formula = "'11.46"
$ws.Range("E26").Value = "  -6.17%  "
$ws.Range("D27").Formula = "'10.60"
$ws.Range("E27").Value = "  -4.36%  "
$ws.Range("D28").Formula = "'36.47"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").Formula = "'13.14"
$ws.Range("E29").Value = "  -2.91%  "
$ws.Range("D30").Formula = "'673.15"
$ws.Range("E30").Value = "  -5.39%  "
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("D33").Formula = "'69.14"
$ws.Range("E33").Value = "  +12.85%  "
$ws.Range("D34").Value = "0.0₃0881"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").Value = "  +11.02%  "
$ws.Range("E36").Value = "  -2.15%  "
$ws.Range("D37").Formula = "'40.06"
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("D38").Formula = "'0.148"
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Formula = "'3.22"
$ws.Range("D42").Formula = "'0.0480"
$ws.Range("E42").Value = "  -3.69%  "
$ws.Range("E43").Value = "  +7.57%  "
$ws.Range("D44").Formula = "'2.80"
$ws.Range("E44").Value = "  -7.18%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Formula = "'3.40"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Formula = "'0.141"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Formula = "'0.000283"
$ws.Range("E47").Value = "  +16.70%  "
$ws.Range("D48").Value = "0.0₆0357"
$ws.Range("E48").Value = "  +10.93%  "
$ws.Range("D49").Formula = "'2.99"
$ws.Range("E49").Value = "  +6.79%  "
$ws.Range("D50").Value = "2.747.65"
$ws.Range("D51").Formula = "'144.37"
$ws.Range("E51").Value = "  -0.06%  "
